$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed from the source data
# (row 26 = "RM 232", row 28 = "SC 92"). Delete the lower-numbered
# row first is not required since we reference by label-derived
# row numbers computed up front, but deleting from the bottom up
# avoids needing to recompute indices.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Restore / clear individual "missing data" cells to match the target state
$ws.Range("C3").Value = 11.2
$ws.Range("E4").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("E9").Value = -6.8
$ws.Range("E10").Value = -6.1
$ws.Range("E17").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("C23").Value = ""

# C32 corresponds to the "SC 193" row after the row deletions above
# (originally row 34); its missing value is now filled in.
$ws.Range("C32").Value = 10.5
